$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.168335437774658
$ws.Range("B1").Value = 2.358563899993896
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.364940404891968
$ws.Range("E1").Value = 1.234638094902039
